$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 925
$ws.Cells.Item(3, 6).Value = 554
$ws.Cells.Item(4, 6).Value = 236
$ws.Cells.Item(5, 6).Value = 15
$ws.Cells.Item(6, 6).Value = 712
$ws.Cells.Item(7, 6).Value = 324
$ws.Cells.Item(9, 6).Value = 124
$ws.Cells.Item(10, 6).Value = 233
$ws.Cells.Item(11, 6).Value = 192
$ws.Cells.Item(12, 6).Value = 4992
$ws.Cells.Item(15, 6).Value = 473
$ws.Cells.Item(17, 6).Value = 526
$ws.Cells.Item(18, 6).Value = 323
$ws.Cells.Item(21, 6).Value = 11
$ws.Cells.Item(22, 6).Value = 696
$ws.Cells.Item(23, 6).Value = 81
$ws.Cells.Item(24, 6).Value = 293
$ws.Cells.Item(25, 6).Value = 990
$ws.Cells.Item(27, 6).Value = 1692
$ws.Cells.Item(28, 6).Value = 413

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 656
$ws.Cells.Item(5, 6).Value = 258
$ws.Cells.Item(6, 6).Value = 40
$ws.Cells.Item(7, 6).Value = 252

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 163

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 925
$ws.Cells.Item(6, 6).Value = 163
$ws.Cells.Item(7, 6).Value = 554
$ws.Cells.Item(8, 6).Value = 236
$ws.Cells.Item(9, 6).Value = 16
$ws.Cells.Item(10, 6).Value = 712
$ws.Cells.Item(11, 6).Value = 656
$ws.Cells.Item(12, 6).Value = 324
$ws.Cells.Item(14, 6).Value = 124
$ws.Cells.Item(15, 6).Value = 233
$ws.Cells.Item(16, 6).Value = 192
$ws.Cells.Item(17, 6).Value = 4992
$ws.Cells.Item(20, 6).Value = 258
$ws.Cells.Item(21, 6).Value = 473
$ws.Cells.Item(23, 6).Value = 526
$ws.Cells.Item(24, 6).Value = 323
$ws.Cells.Item(26, 6).Value = 40
$ws.Cells.Item(28, 6).Value = 11
$ws.Cells.Item(29, 6).Value = 252
$ws.Cells.Item(32, 6).Value = 696
$ws.Cells.Item(36, 6).Value = 81
$ws.Cells.Item(37, 6).Value = 293
$ws.Cells.Item(38, 6).Value = 990
$ws.Cells.Item(40, 6).Value = 1692
$ws.Cells.Item(41, 6).Value = 413
